$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds plain-text numbers (dotted thousands
# separators like "66.611.28", or decimals like "609.61"). Force the
# number format to Text before assigning so Excel does not silently
# convert/reformat values that look numeric (e.g. "25.50" -> 25.5).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.611.28'
$ws.Range("E2").Value = '  +0.94%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.602.78'

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.61'
$ws.Range("E5").Value = '  +0.75%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.21'
$ws.Range("E6").Value = '  +3.58%  '

$ws.Range("E7").Value = '  +0.13%  '

$ws.Range("E8").Value = '  -0.58%  '

$ws.Range("E9").Value = '  +1.99%  '

$ws.Range("E10").Value = '  +0.19%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.417'
$ws.Range("E11").Value = '  +1.11%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.215.37'
$ws.Range("E12").Value = '  +1.65%  '

$ws.Range("E13").Value = '  +1.64%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '29.84'
$ws.Range("E14").Value = '  -0.36%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.604.03'
$ws.Range("E15").Value = '  +1.60%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.706.81'
$ws.Range("E16").Value = '  +0.94%  '

$ws.Range("E17").Value = '  +0.84%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.57'
$ws.Range("E18").Value = '  +2.53%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.39'
$ws.Range("E19").Value = '  +3.52%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.14'
$ws.Range("E20").Value = '  +2.17%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '428.68'
$ws.Range("E21").Value = '  -0.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.620'
$ws.Range("E22").Value = '  +1.67%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '78.91'
$ws.Range("E23").Value = '  -0.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.752.06'
$ws.Range("E24").Value = '  +1.74%  '

$ws.Range("E25").Value = '  +0.01%  '

$ws.Range("E26").Value = '  +5.32%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.33'
$ws.Range("E27").Value = '  +5.10%  '

$ws.Range("E28").Value = '  +4.90%  '

$ws.Range("E29").Value = '  +0.53%  '

$ws.Range("E30").Value = '  -0.08%  '

$ws.Range("E31").Value = '  +1.40%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.601.32'
$ws.Range("E32").Value = '  +1.69%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.159'
$ws.Range("E33").Value = '  +3.65%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.50'
$ws.Range("E34").Value = '  +0.02%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.88'
$ws.Range("E35").Value = '  +0.32%  '

$ws.Range("E36").Value = '  +0.00%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.68'
$ws.Range("E37").Value = '  +1.49%  '

$ws.Range("E38").Value = '  -1.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '177.86'
$ws.Range("E39").Value = '  +2.42%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0858'
$ws.Range("E40").Value = '  +1.08%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.26'
$ws.Range("E41").Value = '  +1.28%  '

$ws.Range("E42").Value = '  +1.04%  '

$ws.Range("E43").Value = '  -0.06%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.57'
$ws.Range("E44").Value = '  +10.61%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.13'
$ws.Range("E46").Value = '  -1.06%  '

$ws.Range("E47").Value = '  -1.43%  '

$ws.Range("E48").Value = '  +3.05%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.20'
$ws.Range("E49").Value = '  +1.36%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.956'
$ws.Range("E50").Value = '  +1.97%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.432.15'
$ws.Range("E51").Value = '  +5.81%  '
